$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.055.10'
$ws.Range('E2').Value = '  +0.47%  '
$ws.Range('D3').Value = '1.680.41'
$ws.Range('E3').Value = '  +0.84%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '215.90'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.17%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.518'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -2.78%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  +1.93%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '21.39'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +5.72%  '
$ws.Range('E10').Value = '  +0.58%  '
$ws.Range('E11').Value = '  -0.90%  '
$ws.Range('D12').Value = '1.917.88'
$ws.Range('E12').Value = '  +0.86%  '
$ws.Range('D13').Value = '1.709.05'
$ws.Range('E13').Value = '  +2.28%  '
$ws.Range('E14').Value = '  +0.86%  '
$ws.Range('E15').Value = '  +1.66%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '66.29'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.30%  '
$ws.Range('D17').Value = '27.041.40'
$ws.Range('E17').Value = '  +0.43%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '8.16'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +1.69%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '236.08'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.34%  '
$ws.Range('D20').Value = '0.0₃0737'
$ws.Range('E20').Value = '  +0.69%  '
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('E22').Value = '  +2.98%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.26'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.88%  '
$ws.Range('E24').Value = '  -3.39%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '147.03'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.74%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.27'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +2.11%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '16.49'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +3.73%  '
$ws.Range('E28').Value = '  -1.65%  '
$ws.Range('E29').Value = '  +0.12%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0497'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.06%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.18'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.07%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.37'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.17%  '
$ws.Range('D33').Value = '1.544.22'
$ws.Range('E33').Value = '  +6.13%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.19'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +1.41%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.72'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +5.49%  '
$ws.Range('E36').Value = '  -0.75%  '
$ws.Range('E37').Value = '  +1.76%  '
$ws.Range('E38').Value = '  +1.31%  '
$ws.Range('E39').Value = '  +3.11%  '
$ws.Range('E40').Value = '  +6.77%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '67.86'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +3.13%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.54'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.91%  '
$ws.Range('E44').Value = '  -0.78%  '
$ws.Range('D45').Value = '1.822.63'
$ws.Range('E46').Value = '  -0.48%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '90.56'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.03%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0108'
$ws.Range('E48').Value = '  +3.19%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.54'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.45%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.104'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.82%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '8.02'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +5.98%  '
